$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Rows that previously had no "Sprint Number" (column E) entry - copy the
# existing column-E formatting (style index 2, matching the rest of the
# table) from a populated neighbor before writing the new value.
$newSprintRows = @{
    8  = 3
    9  = 3
    12 = 3
    14 = 4
    15 = 3
    16 = 4
    19 = 3
    24 = 3
    27 = 3
}

$ws.Range("E10").Copy() | Out-Null
foreach ($row in $newSprintRows.Keys) {
    $cell = $ws.Range("E$row")
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $cell.Value = $newSprintRows[$row]
}
$excel.CutCopyMode = 0

# Rows 20 and 21 already had a Sprint Number; bump them to reflect the
# newly added sprint 4 artifacts.
$ws.Range("E20").Value = 4
$ws.Range("E21").Value = 4
